$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 324 (existing rows 324-358 shift down to 326-360).
$ws.Rows.Item(324).Insert()
$ws.Rows.Item(324).Insert()

# Populate the first new row (324) - weekly "Primera" entry for 2021-10-06 (serial 44491)
$ws.Range("A324").Value = 9
$ws.Range("B324").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C324").Value = "Metropolitana"
$ws.Range("D324").Value = 44491
$ws.Range("E324").Value = 13
$ws.Range("F324").Value = 100112009
$ws.Range("G324").Value = "Acelga"
$ws.Range("H324").Value = "Sin especificar"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 70
$ws.Range("K324").Value = 11000
$ws.Range("L324").Value = 12000
$ws.Range("M324").Value = 11500
$ws.Range("N324").Value = "$/docena de atados"
$ws.Range("O324").Value = "Región Metropolitana"
$ws.Range("P324").Value = 3833
$ws.Range("Q324").Value = 3
$ws.Range("R324").Value = "Hortaliza"

# Populate the second new row (325) - weekly "Segunda" entry for 2021-10-06 (serial 44491)
$ws.Range("A325").Value = 9
$ws.Range("B325").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C325").Value = "Metropolitana"
$ws.Range("D325").Value = 44491
$ws.Range("E325").Value = 13
$ws.Range("F325").Value = 100112009
$ws.Range("G325").Value = "Acelga"
$ws.Range("H325").Value = "Sin especificar"
$ws.Range("I325").Value = "Segunda"
$ws.Range("J325").Value = 43
$ws.Range("K325").Value = 9000
$ws.Range("L325").Value = 10000
$ws.Range("M325").Value = 9512
$ws.Range("N325").Value = "$/docena de atados"
$ws.Range("O325").Value = "Región Metropolitana"
$ws.Range("P325").Value = 3171
$ws.Range("Q325").Value = 3
$ws.Range("R325").Value = "Hortaliza"
